# InputParam1.xlsx — introduce junctional/bulk compartment parameters.
#
# B_SRtot / K_SRBuffer are renamed in place to B_CSQ / K_CSQ (same values),
# and eight new diffusion-related parameters are appended below the
# existing parameter table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename existing SR-buffer parameters to CSQ naming ---------------
$ws.Range("A96").Value = "B_CSQ"
$ws.Range("A97").Value = "K_CSQ"

# --- Clear the leftover pre-formatted (empty) cell at B103 so that the
#     scientific-notation number format can be relocated to B104, where
#     the new D_Pi value belongs. -------------------------------------
$ws.Range("B103").ClearFormats()

# --- New parameter rows -------------------------------------------------
# (entered in the order that reproduces the shared-string table ordering
#  of the target workbook: rows 99-105 first, row 98 last)
$ws.Range("A99").Value = "diffusion_length"
$ws.Range("B99").Value = 0.1

$ws.Range("A100").Value = "D_ion"
$ws.Range("B100").Value = 1000

$ws.Range("A101").Value = "D_ATP"
$ws.Range("B101").Value = 500

$ws.Range("A102").Value = "D_parv"
$ws.Range("B102").Value = 1

$ws.Range("A103").Value = "D_CSQ"
$ws.Range("B103").Value = 0.1

$ws.Range("A104").Value = "D_Pi"
$ws.Range("B104").Value = 500
$ws.Range("B104").NumberFormat = "0.000000000"

$ws.Range("A105").Value = "tau_V"
$ws.Range("B105").Value = 0.00015

$ws.Range("A98").Value = "k_offCSQ"
$ws.Range("B98").Value = 5

# --- Cosmetic view-state updates (best effort) --------------------------
$ws.Columns.Item(2).ColumnWidth = 12.9

$win = $wb.Windows.Item(1)
$win.ScrollRow = 82
$win.ScrollColumn = 1

$ws.Range("E99").Select()
